{"js": "// The author shortened \"Philip\" to \"Phil\" in the byline \"By Professor\n// Philip Hubbard\" (i.e. deleted the \"ip\" out of \"Philip\"). Word marks\n// the spot of the most recent edit with its hidden \"_GoBack\" bookmark,\n// so that bookmark moves from its old location (near the word\n// \"Outline\", further down the document) to the new edit point, right\n// between \"Phil\" and \" Hubbard\".\n\n// 1. Drop the stale _GoBack bookmark wherever it currently lives.\nconst existing = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\nexisting.load(\"isNullObject\");\nawait context.sync();\nif (!existing.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// 2. Shrink \"Philip\" down to \"Phil\" inside the byline by deleting the\n//    \"ip\" characters (scoped to the \"Philip\" match so the other,\n//    unrelated \"ip\" elsewhere in the document is left untouched).\nconst nameMatches = context.document.body.search(\"Philip\", { matchCase: true, matchWholeWord: false });\nnameMatches.load(\"items\");\nawait context.sync();\nconst philipRange = nameMatches.items[0];\n\nconst toRemove = philipRange.search(\"ip\", { matchCase: true, matchWholeWord: false });\ntoRemove.load(\"items\");\nawait context.sync();\ntoRemove.items[0].delete();\nawait context.sync();\n\n// 3. Re-locate \"Phil\" and plant a fresh, empty _GoBack bookmark right\n//    after it (i.e. exactly between \"Phil\" and \" Hubbard\"), mirroring\n//    where Word leaves the cursor after such an edit.\nconst editMatches = context.document.body.search(\"Phil\", { matchCase: true, matchWholeWord: false });\neditMatches.load(\"items\");\nawait context.sync();\nconst editPoint = editMatches.items[0].getRange(Word.RangeLocation.end);\neditPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The author shortened \"Philip\" to \"Phil\" in the byline \"By Professor\n# Philip Hubbard\" (i.e. deleted the \"ip\" out of \"Philip\"). Word marks\n# the spot of the most recent edit with its hidden \"_GoBack\" bookmark,\n# so that bookmark moves from its old location (near the word\n# \"Outline\", further down the document) to the new edit point, right\n# between \"Phil\" and \" Hubbard\".\n\n$d = $word.ActiveDocument\n\n# 1. Drop the stale _GoBack bookmark wherever it currently lives.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2. Shrink \"Philip\" down to \"Phil\" inside the byline.\n$find = $d.Content.Find\n$find.Text = \"Philip\"\n$find.Replacement.Text = \"Phil\"\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, \"Phil\", 2)\n\n# 3. Re-locate \"Phil\" and plant a fresh, empty _GoBack bookmark right\n#    after it (i.e. exactly between \"Phil\" and \" Hubbard\"), mirroring\n#    where Word leaves the cursor after such an edit.\n$editPoint = $d.Content\n$editPoint.Find.Text = \"Phil\"\n$editPoint.Find.MatchCase = $true\n$editPoint.Find.Execute()\n$editPoint.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $editPoint)\n"}
